# [Fonds de solidarite] Add 2020-12-09 data
# Refreshes the "nombre_aides" (col C) and "montant_total" (col D) figures
# for a handful of rows. Source cells are stored as text, so we force the
# Text number format before writing to keep the values as strings (avoids
# Excel auto-converting the numeric-looking text into real numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 4  - Auvergne-Rhone-Alpes / SARL
Set-TextValue "C4" "1437"
Set-TextValue "D4" "9595127.25"

# Row 6  - Auvergne-Rhone-Alpes
Set-TextValue "C6" "1020"
Set-TextValue "D6" "5886590.53"

# Row 24 - Centre-Val de Loire
Set-TextValue "C24" "202"
Set-TextValue "D24" "1090826.23"

# Row 51 - Hauts-de-France / SARL
Set-TextValue "C51" "1143"
Set-TextValue "D51" "8291040.22"

# Row 52 - Hauts-de-France
Set-TextValue "C52" "794"
Set-TextValue "D52" "5030000.28"

# Row 55 - Ile-de-France / Entrepreneur individuel
Set-TextValue "C55" "10144"
Set-TextValue "D55" "29325153.25"

# Row 61 - Ile-de-France
Set-TextValue "C61" "6684"
Set-TextValue "D61" "28881781.83"

# Row 96 - Pays de la Loire / SARL
Set-TextValue "C96" "655"
Set-TextValue "D96" "4458316.04"

# Row 101 - Provence-Alpes-Cote d'Azur
Set-TextValue "C101" "1497"
Set-TextValue "D101" "3795586.09"
